$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; everything from old B..F shifts to C..G
$ws.Columns.Item(2).Insert()

# New header "Date" in B1 (same style as the rest of the header row)
$ws.Cells.Item(1, 2).Value = "Date"

# New date value in B2, styled as center/top aligned with a short-date format
$ws.Cells.Item(2, 2).HorizontalAlignment = -4108
$ws.Cells.Item(2, 2).VerticalAlignment = -4160
$ws.Cells.Item(2, 2).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2, 2).Value = (Get-Date -Year 2012 -Month 10 -Day 22).Date

# Updated / expanded comment text, now wrapped, in C2
$ws.Cells.Item(2, 3).Value = "il faut que le résultat obtenu soit proportionnel au temps de l'explosion et que cela soit un entier"
$ws.Cells.Item(2, 3).WrapText = $true
$ws.Cells.Item(2, 3).VerticalAlignment = -4160

# Row 2 cells get vertical=top alignment; E2 loses its old center alignment
$ws.Cells.Item(2, 1).VerticalAlignment = -4160
$ws.Cells.Item(2, 4).VerticalAlignment = -4160
$ws.Cells.Item(2, 5).HorizontalAlignment = -4142
$ws.Cells.Item(2, 5).VerticalAlignment = -4160
$ws.Cells.Item(2, 6).VerticalAlignment = -4160
$ws.Cells.Item(2, 7).VerticalAlignment = -4160

# Row 2 custom height
$ws.Rows.Item(2).RowHeight = 31.5

# Extend the bordered block: 9 more rows (21-29), column F only, copying the
# same (border-less, centered) formatting used by the existing E/F filler cells
$ws.Cells.Item(20, 5).Copy()
$ws.Range("F21:F29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C13").Select()
